$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsPerform = $wb.Worksheets.Item("演出")
$wsAll = $wb.Worksheets.Item("全部类型")

# 展览 sheet updates
$wsExhibit.Range("F2").Value = 1601
$wsExhibit.Range("G2").Value = "不可售"
$wsExhibit.Range("F3").Value = 372
$wsExhibit.Range("F4").Value = 5233
$wsExhibit.Range("F6").Value = 10472
$wsExhibit.Range("F7").Value = 267
$wsExhibit.Range("F8").Value = 569
$wsExhibit.Range("F9").Value = 128
$wsExhibit.Range("F10").Value = 142
$wsExhibit.Range("F11").Value = 835

# 演出 sheet updates
$wsPerform.Range("F4").Value = 18

# 全部类型 sheet updates
$wsAll.Range("F2").Value = 1601
$wsAll.Range("G2").Value = "不可售"
$wsAll.Range("F3").Value = 372
$wsAll.Range("F6").Value = 5233
$wsAll.Range("F8").Value = 18
$wsAll.Range("F9").Value = 10472
$wsAll.Range("F10").Value = 267
$wsAll.Range("F11").Value = 569
$wsAll.Range("F12").Value = 128
$wsAll.Range("F15").Value = 142
$wsAll.Range("F16").Value = 835
